$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PPP-3")

$ws.Range("E3").Value = 15
$ws.Range("E4").Value = 11
$ws.Range("E5").Value = 20
$ws.Range("E6").Value = 11
$ws.Range("E7").Value = 6
$ws.Range("E8").Value = 7
$ws.Range("E9").Value = 18
$ws.Range("E10").Value = 18
$ws.Range("E11").Value = 16
$ws.Range("E12").Value = 12
$ws.Range("E13").Value = 17
$ws.Range("E14").Value = 15
$ws.Range("E15").Value = 7
$ws.Range("E16").Value = 10
$ws.Range("E17").Value = 17
$ws.Range("E18").Value = 12
$ws.Range("E19").Value = 16
$ws.Range("E20").Value = 9
$ws.Range("E21").Value = 8
$ws.Range("E22").Value = 13
$ws.Range("E23").Value = 16
$ws.Range("E24").Value = 15
$ws.Range("E25").Value = 20
$ws.Range("E26").Value = 18
$ws.Range("E27").Value = 5
$ws.Range("E29").Value = 17
$ws.Range("E30").Value = 8
$ws.Range("E31").Value = 20
$ws.Range("E32").Value = 20
$ws.Range("E33").Value = 12
$ws.Range("E34").Value = 17
$ws.Range("E35").Value = 10
$ws.Range("E36").Value = 20
$ws.Range("E38").Value = 11
$ws.Range("E39").Value = 11
$ws.Range("E41").Value = 6
$ws.Range("E42").Value = 20
$ws.Range("E43").Value = 20
$ws.Range("E44").Value = 12
$ws.Range("E45").Value = 10
$ws.Range("E46").Value = 15
$ws.Range("E47").Value = 17
$ws.Range("E48").Value = 11
$ws.Range("E49").Value = 9
$ws.Range("E50").Value = 20
$ws.Range("E51").Value = 7
$ws.Range("E53").Value = 7
$ws.Range("E54").Value = 8
$ws.Range("E55").Value = 12
$ws.Range("E56").Value = 8
$ws.Range("E58").Value = 19
$ws.Range("E59").Value = 5
$ws.Range("E60").Value = 7
$ws.Range("E61").Value = 11
$ws.Range("E62").Value = 18
